$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price column cells whose new values would otherwise be
# auto-parsed as numbers (losing exact text formatting like trailing
# zeros) to stay as plain text, matching the source inlineStr cells.
$ws.Range("D4:D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24:D32").NumberFormat = "@"
$ws.Range("D34:D51").NumberFormat = "@"

# Updated coin data scraped on Sat Feb 25 22:27:09 UTC 2023
$ws.Range("D2").Value = "23.049.51"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.581.97"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "299.69"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").Value = "0.3753"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "0.3564"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").Value = "50.40"
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "1.221"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").Value = "0.07981"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "21.98"
$ws.Range("E13").Value = "  -4.06%  "
$ws.Range("D14").Value = "6.423"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "7.309"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "0.00001226"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").Value = "1.580.78"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "92.30"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  -3.19%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "6.347"
$ws.Range("E22").Value = "  -3.19%  "
$ws.Range("D23").Value = "23.037.37"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "12.63"
$ws.Range("E24").Value = "  -3.55%  "
$ws.Range("D25").Value = "2.372"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").Value = "2.848"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "20.62"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "148.28"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("D29").Value = "5.172"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").Value = "131.01"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "2.335"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "6.529"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("D33").Value = "1.753.95"
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "0.9352"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("D35").Value = "0.07347"
$ws.Range("E35").Value = "  -4.56%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "9.974"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "0.08748"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "0.02647"
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("D39").Value = "0.2471"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").Value = "5.995"
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("D41").Value = "1.339"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("D42").Value = "0.6899"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("D43").Value = "11.97"
$ws.Range("E43").Value = "  -6.32%  "
$ws.Range("D44").Value = "14.63"
$ws.Range("E44").Value = "  -7.33%  "
$ws.Range("D45").Value = "0.9993"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "0.6349"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").Value = "3.969"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "2.245"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "130.24"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "0.07862"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "1.188"
$ws.Range("E51").Value = "  +1.60%  "
